$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 previously contained the correct metadata string followed by many
# accidentally-duplicated "Medieval" suffixes; trim it back down.
$ws.Range("B2").Value = "name=Age,dataType=text,updateCriteria=true"

# B3 had the same accidental duplication bug; should just read "Age".
$ws.Range("B3").Value = "Age"

# New subentity row (row 4): id 2, "Medieval" age value, visibilityGroups -1.
$ws.Cells.Item(4, 1).Value = 2
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Cells.Item(4, 2).Value = "Medieval"
$ws.Cells.Item(4, 3).Value = -1
